$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels: A1 becomes "Year", B1 becomes "totalNetAssets"
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "totalNetAssets"

# Update the active selection to B1
$ws.Range("B1").Select()
